$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 519-520, pushing the existing rows
# (old 519..555) down to new rows 521..557. This also makes the
# trailing "new" rows 556/557 fall out naturally from the shift.
$ws.Rows("519:520").Insert()

$newDate = Get-Date -Year 2023 -Month 12 -Day 5 -Hour 0 -Minute 0 -Second 0

# New row 519: Región Metropolitana entry, $/caja 36 atados
$ws.Cells.Item(519, 1).Value = 4
$ws.Cells.Item(519, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(519, 3).Value = "Los Lagos"
$ws.Cells.Item(519, 4).Value = $newDate
$ws.Cells.Item(519, 5).Value = 10
$ws.Cells.Item(519, 6).Value = 100112040
$ws.Cells.Item(519, 7).Value = "Cilantro"
$ws.Cells.Item(519, 8).Value = "Sin especificar"
$ws.Cells.Item(519, 9).Value = "Primera"
$ws.Cells.Item(519, 10).Value = 160
$ws.Cells.Item(519, 11).Value = 24000
$ws.Cells.Item(519, 12).Value = 24000
$ws.Cells.Item(519, 13).Value = 24000
$ws.Cells.Item(519, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(519, 15).Value = "Región Metropolitana"
$ws.Cells.Item(519, 16).Value = 667
$ws.Cells.Item(519, 17).Value = 36
$ws.Cells.Item(519, 18).Value = "Hortaliza"

# New row 520: Región de La Araucanía entry, $/docena de atados (2 kilos)
$ws.Cells.Item(520, 1).Value = 4
$ws.Cells.Item(520, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(520, 3).Value = "Los Lagos"
$ws.Cells.Item(520, 4).Value = $newDate
$ws.Cells.Item(520, 5).Value = 10
$ws.Cells.Item(520, 6).Value = 100112040
$ws.Cells.Item(520, 7).Value = "Cilantro"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 120
$ws.Cells.Item(520, 11).Value = 11000
$ws.Cells.Item(520, 12).Value = 11000
$ws.Cells.Item(520, 13).Value = 11000
$ws.Cells.Item(520, 14).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(520, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(520, 16).Value = 5500
$ws.Cells.Item(520, 17).Value = 2
$ws.Cells.Item(520, 18).Value = "Hortaliza"
